# "data for csc computations"
# Convert the absolute ForestCarbon2022 data paths to relative "data/..." paths
# and add a new "graph" row pointing at the canal network matrix pickle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows: strip the old absolute prefix, keep the trailing relative path.
$ws.Range("B4").Value = "data/dtm_depth_padded.tif"
$ws.Range("B5").Value = "data/depth_extended.tif"
$ws.Range("B6").Value = "data/199_canalblocks_20191008b.shp"
$ws.Range("B7").Value = "data/weather_station_coordinates.xlsx"

# "mesh" row moves up one slot in the shared-string table but keeps its position/content.
$ws.Range("A8").Value = "mesh"
$ws.Range("B8").Value = "data/mesh_0.05.msh2"

# New row for the canal-network graph pickle.
$ws.Range("A9").Value = "graph"
$ws.Range("B9").Value = "data/canal_network_matrix_50meters.p"

# Match the author's final on-screen selection.
$ws.Range("L8").Select()
